$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers (human-readable labels)
$ws.Range("A1").Value = "Superficie útil"
$ws.Range("B1").Value = "Comarca nombre"
$ws.Range("C1").Value = "Número hogares"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Municipio código"
$ws.Range("H1").Value = "Provincia nombre"
$ws.Range("I1").Value = "Municipio nombre"

# Row 2 - dimension/measure identifiers
$ws.Range("A2").Value = "iaest-dimension:superficie-util"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "iaest-measure:numero-hogares"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3 - type (dim / medida)
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "dim"
$ws.Range("I3").Value = "dim"

# Row 4 - data type / codelist references
$ws.Range("A4").Value = "skos:Concept"
$ws.Range("B4").Value = "URI-comarca"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "URI-Provincia"
$ws.Range("I4").Value = "URI-Municipio"

# Row 5 - the mapping file reference moves from column G to column A
$ws.Range("G5").Delete()
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "mapping-superficie-util.xlsx"
